# Applies the commit "Elimna EC anteriores y se agregan nuevos, se modifica
# base de datos" to the EC (Estado de Cuenta) workbook:
#   - Adds a new worker record (MARIA BERNARDA SILVA FUENTES) as the first
#     data row, pushing the existing LUIS ALBERTO PUELLO CASTELLON rows (and
#     the signature block below) down by one row.
#   - Updates the "VALOR MORA" total and the worker/period counters.
#   - Re-flows the signature block (two new rows: the underline and the
#     "NOMBRE/FIRMA DEL REPRESENTANTE LEGAL" captions).
#   - Refreshes the bestFit column widths for the new, wider content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new data row at row 16 ----------------------------------
# This shifts the existing data rows (16-18 -> 17-19) and the signature
# block (23-24 -> 24-25) down by one, carrying their original formatting
# (including the special "last row" bottom border) down with them.
$ws.Rows("16:16").Insert()

# The freshly inserted row 16 picks up a blank/synthesized style; restore
# the normal data-row formatting by pulling it from row 17 (an untouched
# copy of the original row-16 format that got pushed down).
$ws.Range("B17:J17").Copy()
$ws.Range("B16:J16").PasteSpecial(-4122)

# --- 2. Populate the new worker row --------------------------------------
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1047501351"
$ws.Range("D16").Value = "MARIA BERNARDA SILVA FUENTES"
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 60000
$ws.Range("G16").Value = 1500000

# The existing LUIS ALBERTO PUELLO CASTELLON periods are re-ordered
# (descending) while they get pushed down a row.
$ws.Range("E17").Value = "1712"
$ws.Range("E18").Value = "1711"
$ws.Range("E19").Value = "1710"

# --- 3. Update the header totals -----------------------------------------
$ws.Range("E11").Value = 148527
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 4

# --- 4. Fix up the signature block text -----------------------------------
# After the row insert, row 24 already holds the underline text and row 25
# already holds the "NOMBRE/FIRMA" captions (carried down automatically),
# so no text changes are required there - just confirm/re-assert them.
$ws.Range("B24").Value = "___________________________________"
$ws.Range("H24").Value = "___________________________________"
$ws.Range("B25").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H25").Value = "FIRMA DEL REPRESENTANTE LEGAL"

# --- 5. Refresh auto-fit column widths for the new, wider content --------
$ws.Columns("B:J").AutoFit()
